# Cover Letter - Showbie.docx edit script
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# 1. Date: "February 23, 2020" -> "March 6, 2020"
Replace-Text "February 23, 2020" "March 6, 2020"

# 2. Salutation: "Dear Hiring Manager," -> "To Whom It May Concern,"
Replace-Text "Dear Hiring Manager," "To Whom It May Concern,"

# 3. "interested in Showbie's" -> "interested in applying for Showbie's"
Replace-Text "interested in Showbie" "interested in applying for Showbie"

# 4. "creating useful and tested solutions" -> "developing useful and tested solutions"
Replace-Text "creating useful and tested solutions" "developing useful and tested solutions"

# 5. Hackathon paragraph rewrite
Replace-Text "built an ML Algorithm application where we used React Native to create the front-end. For the project, we initially used AWS Amplify and S3, however we moved away from AWS services. " `
             "built an Machine Learning (ML) Algorithm application where we used React Native to create the front-end. For the hackathon, we initially used AWS Amplify and S3, however we moved away from AWS services towards the end. "

# 6. "I would bring to any position" -> "I would bring to this position"
Replace-Text "I would bring to any position" "I would bring to this position"

# 7. github line rewrite
Replace-Text " full of my projects that details my projects and experiences in software engineering. I can be reached anytime via " `
             " that details my projects and experiences in software engineering and full-stack development. I can be reached anytime via "
